$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 4).Value = 44573
$ws.Cells.Item(4, 13).Value = 300
$ws.Cells.Item(4, 14).Value = 20500
$ws.Cells.Item(4, 15).Value = 21000
$ws.Cells.Item(4, 16).Value = 20750
$ws.Cells.Item(4, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(4, 19).Value = 1153
$ws.Cells.Item(5, 4).Value = 44573
$ws.Cells.Item(5, 13).Value = 400
$ws.Cells.Item(5, 14).Value = 17500
$ws.Cells.Item(5, 15).Value = 18000
$ws.Cells.Item(5, 16).Value = 17750
$ws.Cells.Item(5, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(5, 19).Value = 986
$ws.Cells.Item(6, 4).Value = 44566
$ws.Cells.Item(6, 11).Value = 'Modesto'
$ws.Cells.Item(6, 13).Value = 100
$ws.Cells.Item(6, 14).Value = 23000
$ws.Cells.Item(6, 15).Value = 24000
$ws.Cells.Item(6, 16).Value = 23500
$ws.Cells.Item(6, 19).Value = 1306
$ws.Cells.Item(7, 4).Value = 44566
$ws.Cells.Item(7, 11).Value = 'Modesto'
$ws.Cells.Item(7, 13).Value = 160
$ws.Cells.Item(7, 14).Value = 21000
$ws.Cells.Item(7, 15).Value = 22000
$ws.Cells.Item(7, 16).Value = 21500
$ws.Cells.Item(7, 19).Value = 1194
$ws.Cells.Item(8, 12).Value = 'Especial'
$ws.Cells.Item(8, 13).Value = 340
$ws.Cells.Item(8, 14).Value = 22500
$ws.Cells.Item(8, 15).Value = 23000
$ws.Cells.Item(8, 16).Value = 22750
$ws.Cells.Item(8, 19).Value = 1264
$ws.Cells.Item(9, 4).Value = 44545
$ws.Cells.Item(9, 11).Value = 'Castle Brite'
$ws.Cells.Item(9, 12).Value = 'Primera'
$ws.Cells.Item(9, 13).Value = 400
$ws.Cells.Item(9, 14).Value = 20500
$ws.Cells.Item(9, 15).Value = 21000
$ws.Cells.Item(9, 16).Value = 20750
$ws.Cells.Item(9, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(9, 19).Value = 1153
$ws.Cells.Item(10, 4).Value = 44545
$ws.Cells.Item(10, 11).Value = 'Castle Brite'
$ws.Cells.Item(10, 12).Value = 'Segunda'
$ws.Cells.Item(10, 13).Value = 300
$ws.Cells.Item(10, 14).Value = 15500
$ws.Cells.Item(10, 15).Value = 16000
$ws.Cells.Item(10, 16).Value = 15750
$ws.Cells.Item(10, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(10, 19).Value = 875
$ws.Cells.Item(11, 4).Value = 44580
$ws.Cells.Item(11, 11).Value = 'Modesto'
$ws.Cells.Item(11, 13).Value = 300
$ws.Cells.Item(11, 14).Value = 22500
$ws.Cells.Item(11, 15).Value = 23000
$ws.Cells.Item(11, 16).Value = 22750
$ws.Cells.Item(11, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(11, 19).Value = 1264
$ws.Cells.Item(12, 4).Value = 44580
$ws.Cells.Item(12, 11).Value = 'Modesto'
$ws.Cells.Item(12, 13).Value = 400
$ws.Cells.Item(12, 14).Value = 19500
$ws.Cells.Item(12, 15).Value = 20000
$ws.Cells.Item(12, 16).Value = 19750
$ws.Cells.Item(12, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(12, 19).Value = 1097
$ws.Cells.Item(13, 4).Value = 44160
$ws.Cells.Item(13, 12).Value = 'Primera'
$ws.Cells.Item(13, 13).Value = 240
$ws.Cells.Item(13, 14).Value = 20500
$ws.Cells.Item(13, 15).Value = 21000
$ws.Cells.Item(13, 16).Value = 20750
$ws.Cells.Item(13, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(13, 19).Value = 1383
$ws.Cells.Item(13, 20).Value = 15
$ws.Cells.Item(14, 4).Value = 44552
$ws.Cells.Item(14, 12).Value = 'Especial'
$ws.Cells.Item(14, 13).Value = 360
$ws.Cells.Item(14, 14).Value = 20000
$ws.Cells.Item(14, 16).Value = 20500
$ws.Cells.Item(14, 19).Value = 1139
$ws.Cells.Item(15, 4).Value = 44552
$ws.Cells.Item(15, 13).Value = 280
$ws.Cells.Item(15, 14).Value = 18000
$ws.Cells.Item(15, 15).Value = 19000
$ws.Cells.Item(15, 16).Value = 18500
$ws.Cells.Item(15, 19).Value = 1028
$ws.Cells.Item(16, 4).Value = 44175
$ws.Cells.Item(16, 11).Value = 'Castle Brite'
$ws.Cells.Item(16, 12).Value = 'Primera'
$ws.Cells.Item(16, 13).Value = 300
$ws.Cells.Item(16, 14).Value = 21000
$ws.Cells.Item(16, 15).Value = 22000
$ws.Cells.Item(16, 16).Value = 21500
$ws.Cells.Item(16, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(16, 19).Value = 1194
$ws.Cells.Item(17, 4).Value = 44553
$ws.Cells.Item(17, 12).Value = 'Especial'
$ws.Cells.Item(17, 13).Value = 360
$ws.Cells.Item(17, 14).Value = 23000
$ws.Cells.Item(17, 15).Value = 24000
$ws.Cells.Item(17, 16).Value = 23500
$ws.Cells.Item(17, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(17, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(17, 19).Value = 1469
$ws.Cells.Item(17, 20).Value = 16
$ws.Cells.Item(18, 4).Value = 44553
$ws.Cells.Item(18, 12).Value = 'Primera'
$ws.Cells.Item(18, 14).Value = 21000
$ws.Cells.Item(18, 15).Value = 22000
$ws.Cells.Item(18, 16).Value = 21500
$ws.Cells.Item(18, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(18, 19).Value = 1344
$ws.Cells.Item(18, 20).Value = 16
$ws.Cells.Item(19, 4).Value = 44553
$ws.Cells.Item(19, 12).Value = 'Segunda'
$ws.Cells.Item(19, 13).Value = 240
$ws.Cells.Item(19, 14).Value = 17000
$ws.Cells.Item(19, 16).Value = 17500
$ws.Cells.Item(19, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(19, 19).Value = 1094
$ws.Cells.Item(19, 20).Value = 16
$ws.Cells.Item(20, 4).Value = 44546
$ws.Cells.Item(20, 11).Value = 'Castle Brite'
$ws.Cells.Item(20, 13).Value = 300
$ws.Cells.Item(20, 14).Value = 22500
$ws.Cells.Item(20, 15).Value = 23000
$ws.Cells.Item(20, 16).Value = 22750
$ws.Cells.Item(20, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(20, 19).Value = 1264
$ws.Cells.Item(20, 20).Value = 18
$ws.Cells.Item(21, 4).Value = 44546
$ws.Cells.Item(21, 11).Value = 'Castle Brite'
$ws.Cells.Item(21, 14).Value = 20500
$ws.Cells.Item(21, 15).Value = 21000
$ws.Cells.Item(21, 16).Value = 20750
$ws.Cells.Item(21, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(21, 19).Value = 1153
$ws.Cells.Item(21, 20).Value = 18
$ws.Cells.Item(22, 4).Value = 44189
$ws.Cells.Item(22, 11).Value = 'Dina'
$ws.Cells.Item(22, 12).Value = 'Especial'
$ws.Cells.Item(22, 13).Value = 120
$ws.Cells.Item(22, 14).Value = 23500
$ws.Cells.Item(22, 15).Value = 24000
$ws.Cells.Item(22, 16).Value = 23750
$ws.Cells.Item(22, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(22, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(22, 19).Value = 1319
$ws.Cells.Item(22, 20).Value = 18
$ws.Cells.Item(23, 4).Value = 44189
$ws.Cells.Item(23, 11).Value = 'Dina'
$ws.Cells.Item(23, 12).Value = 'Primera'
$ws.Cells.Item(23, 13).Value = 200
$ws.Cells.Item(23, 14).Value = 21500
$ws.Cells.Item(23, 15).Value = 22000
$ws.Cells.Item(23, 16).Value = 21750
$ws.Cells.Item(23, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(23, 19).Value = 1208
$ws.Cells.Item(24, 4).Value = 44559
$ws.Cells.Item(24, 11).Value = 'Modesto'
$ws.Cells.Item(24, 12).Value = 'Especial'
$ws.Cells.Item(24, 13).Value = 400
$ws.Cells.Item(24, 14).Value = 25000
$ws.Cells.Item(24, 15).Value = 26000
$ws.Cells.Item(24, 16).Value = 25500
$ws.Cells.Item(24, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(24, 19).Value = 1417
$ws.Cells.Item(25, 4).Value = 44559
$ws.Cells.Item(25, 11).Value = 'Modesto'
$ws.Cells.Item(25, 13).Value = 320
$ws.Cells.Item(25, 14).Value = 22000
$ws.Cells.Item(25, 15).Value = 23000
$ws.Cells.Item(25, 16).Value = 22500
$ws.Cells.Item(25, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(25, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(25, 19).Value = 1250
$ws.Cells.Item(25, 20).Value = 18
